# Generate Report for Handoff
# Updates the status from "In Translation" to "Ready for handoff" and
# refreshes the handoff timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-24-12 02:24:58"

# zh-cn detail sheet: ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-12 02:24:55"

# de-de detail sheet: ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-12 02:24:58"
